$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L167:Q167").Formula = "=B167*$I167"
